# Auto-generated edit: refresh cryptos price / Volume(1h) figures to match
# the latest GitHub Actions scrape, and swap the Hedera / PancakeSwap rows
# (rows 32-33) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.017.93"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "3.322.22"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'229.42"
$ws.Range("E5").Value = "  -5.21%  "
$ws.Range("D6").Value = "'616.50"
$ws.Range("E6").Value = "  -4.34%  "
$ws.Range("D7").Value = "'1.35"
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("E8").Value = "  -7.12%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.914"
$ws.Range("E10").Value = "  -8.67%  "
$ws.Range("D11").Value = "3.317.87"
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").Value = "'41.40"
$ws.Range("E12").Value = "  -5.40%  "
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "93.010.23"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "3.951.32"
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").Value = "'7.91"
$ws.Range("E18").Value = "  -8.02%  "
$ws.Range("D19").Value = "3.326.14"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").Value = "'17.06"
$ws.Range("E20").Value = "  -7.11%  "
$ws.Range("D21").Value = "'10.91"
$ws.Range("E21").Value = "  -8.36%  "
$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("D23").Value = "'487.66"
$ws.Range("E23").Value = "  -4.71%  "
$ws.Range("D24").Value = "'0.445"
$ws.Range("E24").Value = "  -10.42%  "
$ws.Range("E25").Value = "  -7.42%  "
$ws.Range("D26").Value = "'6.00"
$ws.Range("E26").Value = "  -10.40%  "
$ws.Range("D27").Value = "'89.20"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "3.504.29"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").Value = "'11.48"
$ws.Range("E29").Value = "  -6.46%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'10.95"
$ws.Range("E31").Value = "  -8.21%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.63"
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.134"
$ws.Range("E33").Value = "  -3.77%  "
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  -7.58%  "
$ws.Range("D36").Value = "'28.08"
$ws.Range("E36").Value = "  -8.76%  "
$ws.Range("D37").Value = "'0.522"
$ws.Range("E37").Value = "  -10.22%  "
$ws.Range("D38").Value = "'522.36"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'7.26"
$ws.Range("E40").Value = "  -7.16%  "
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("E42").Value = "  -8.22%  "
$ws.Range("D43").Value = "'0.863"
$ws.Range("E43").Value = "  -6.02%  "
$ws.Range("D44").Value = "'23.99"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("D47").Value = "'0.0402"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("D48").Value = "'5.27"
$ws.Range("E48").Value = "  -5.26%  "
$ws.Range("D49").Value = "'52.81"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("D51").Value = "'7.78"
$ws.Range("E51").Value = "  -6.23%  "
